$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2"=1.424036728313411; "C2"=0.3236541775035846; "D2"=0.04507898464445503; "E2"=0.06511262675786256; "F2"=1.639132484668778; "H2"=0.07973214163530429; "M2"=0.4641035925004218; "N2"=1.535336609826572
    "B3"=1.285519530909482; "C3"=0.2837497512571474; "D3"=0.04494209123561888; "E3"=0.06121158810409; "F3"=1.574188105213693; "H3"=0.07973214163530429; "M3"=0.4201255521913936; "N3"=1.541218807427242
    "B4"=1.201211014342164; "C4"=0.2593407130880507; "D4"=0.04488412776232309; "E4"=0.05886348798735952; "F4"=1.535440081973988; "H4"=0.07973214163530429; "M4"=0.3933849671949758; "N4"=1.545400875321207
    "B5"=1.167037541836521; "C5"=0.2494159201641253; "D5"=0.04486699105726544; "E5"=0.05791825864116973; "F5"=1.519930377363167; "H5"=0.07973214163530429; "M5"=0.3825524560193969; "N5"=1.547247549127789
    "B6"=1.161374008880784; "C6"=0.247769216471255; "D6"=0.04486453490815734; "E6"=0.05776200080713778; "F6"=1.517371846065316; "H6"=0.07973214163530429; "M6"=0.3807575801321619; "N6"=1.547562762902587
    "B7"=1.200749402902773; "C7"=0.2592067759426584; "D7"=0.04488387049749676; "E7"=0.05885069344601845; "F7"=1.535229781470861; "H7"=0.07973214163530429; "M7"=0.3932386168815967; "N7"=1.545425204728346
    "B8"=1.376119883514548; "C8"=0.309875232842586; "D8"=0.04502632242814997; "E8"=0.06375762959566345; "F8"=1.616503545929476; "H8"=0.07973214163530429; "M8"=0.4488846923908909; "N8"=1.537245818491783
    "B9"=1.726085186550506; "C9"=0.4100269204001847; "D9"=0.04551621620696267; "E9"=0.07376463389199017; "F9"=1.784988902980047; "H9"=0.07973214163530429; "M9"=0.5601563451744198; "N9"=1.525776409137009
    "B10"=1.987179816297555; "C10"=0.4841779654802281; "D10"=0.04600951388986374; "E10"=0.08136644222486211; "F10"=1.914560368363595; "H10"=0.07973214163530429; "M10"=0.6433245139769639; "N10"=1.520198712223134
    "B11"=2.106884454330782; "C11"=0.5180544894154764; "D11"=0.04626399724158148; "E11"=0.08488230762633719; "F11"=1.974813923900314; "H11"=0.07973214163530429; "M11"=0.6814913734584849; "N11"=1.51829365353062
    "B12"=2.152352063190108; "C12"=0.530904948714749; "D12"=0.04636478001035016; "E12"=0.08622224657724331; "F12"=1.99782288424629; "H12"=0.07973214163530429; "M12"=0.6959938909048304; "N12"=1.517664320210883
    "B13"=2.142553593153878; "C13"=0.5281363656361577; "D13"=0.04634287680993765; "E13"=0.08593328213708418; "F13"=1.992858882369205; "H13"=0.07973214163530429; "M13"=0.6928682838775302; "N13"=1.517795744980944
    "B14"=2.11062231097776; "C14"=0.5191112527257928; "D14"=0.04627219968740803; "E14"=0.08499237226350687; "F14"=1.976703010314253; "H14"=0.07973214163530429; "M14"=0.6826835016323685; "N14"=1.518240026001337
    "B15"=2.091081582217896; "C15"=0.5135860355192108; "D15"=0.04622948562904838; "E15"=0.08441715968761798; "F15"=1.96683222068566; "H15"=0.07973214163530429; "M15"=0.6764515312012946; "N15"=1.518524185800914
    "B16"=1.979375929424407; "C16"=0.4819670823924866; "D16"=0.04599349574513667; "E16"=0.08113785584105671; "F16"=1.910649353383747; "H16"=0.07973214163530429; "M16"=0.6408370647884993; "N16"=1.52033603317679
    "B17"=1.911089589105586; "C17"=0.4626079122073179; "D17"=0.04585649130051905; "E17"=0.07914108093479655; "F17"=1.876521370577109; "H17"=0.07973214163530429; "M17"=0.6190752488235631; "N17"=1.521610295426342
    "B18"=1.871900506117413; "C18"=0.4514865818217686; "D18"=0.04578051574333131; "E18"=0.07799801065649348; "F18"=1.857015012658849; "H18"=0.07973214163530429; "M18"=0.6065896944467255; "N18"=1.522402659140653
    "B19"=1.858646645007411; "C19"=0.4477233780577876; "D19"=0.04575527433536308; "E19"=0.07761191034518333; "F19"=1.850431533835945; "H19"=0.07973214163530429; "M19"=0.60236762099629; "N19"=1.522681114199315
    "B20"=1.918349709972404; "C20"=0.4646673167827657; "D20"=0.04587078255372035; "E20"=0.07935307809233905; "F20"=1.880141579567436; "H20"=0.07973214163530429; "M20"=0.6213885835580868; "N20"=1.52146848774494
    "B21"=2.119997522321114; "C21"=0.5217615351793938; "D21"=0.04629283871241796; "E21"=0.08526850604737746; "F21"=1.98144313385751; "H21"=0.07973214163530429; "M21"=0.6856736609248486; "N21"=1.518107021846191
    "B22"=2.252592724921442; "C22"=0.5592056774824528; "D22"=0.0465944562995162; "E22"=0.089184576562225; "F22"=2.048771900394712; "H22"=0.07973214163530429; "M22"=0.7279773082383656; "N22"=1.516447316777899
    "B23"=2.181748969656098; "C23"=0.539208714908284; "D23"=0.04643108834294196; "E23"=0.08708983740776688; "F23"=2.012733320925349; "H23"=0.07973214163530429; "M23"=0.7053720051119825; "N23"=1.517283593028381
    "B24"=1.915067192714616; "C24"=0.4637362335532202; "D24"=0.04586431279417269; "E24"=0.07925721892690518; "F24"=1.878504527427594; "H24"=0.07973214163530429; "M24"=0.6203426452986349; "N24"=1.521532412904023
    "B25"=1.630732611316546; "C25"=0.3828396058991643; "D25"=0.04536063297808823; "E25"=0.07101474698387733; "F25"=1.738410230205346; "H25"=0.07973214163530429; "M25"=0.5298133770013322; "N25"=1.528383477871657
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
